$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("F3").ClearContents()
$ws.Range("K3").Value = "30,46 TL - 60,94 TL - 609,43 TL"

# Row 4
$ws.Range("F4").ClearContents()
$ws.Range("K4").Value = "30,46 TL - 60,94 TL - 609,43 TL"

# Row 5
$ws.Range("F5").ClearContents()
$ws.Range("K5").Value = "30,46 TL - 60,94 TL - 609,43 TL"

# Row 6
$ws.Range("K6").Value = "6,09 TL - 12,19 TL - 152,35 TL"

# Row 8
$ws.Range("F8").ClearContents()
$ws.Range("K8").Value = "15,23 TL - 30,47 TL - 304,71 TL"

# Row 9
$ws.Range("F9").ClearContents()
$ws.Range("K9").Value = "15,23 TL - 30,47 TL - 304,71 TL"

# Row 10
$ws.Range("F10").ClearContents()
$ws.Range("K10").Value = "15,23 TL - 30,47 TL - 304,71 TL"

# Row 11
$ws.Range("K11").Value = "3,05 TL - 6,09 TL - 76,17 TL"

# Row 12
$ws.Range("K12").Value = "WU: ,USD–; Diğer: 404,16 TL–3.403,42 TL"

# Row 13
$ws.Range("D13").Value = "Hesaba: Asgari 1 TL | Azami 300 TL"
$ws.Range("E13").Value = "Hesaba: Asgari 1 TL | Azami 851,5 TL"
$ws.Range("F13").ClearContents()
$ws.Range("K13").Value = "Hesaba: Asgari 1 TL | Azami 865,75 TL"

# Row 14
$ws.Range("F14").ClearContents()
$ws.Range("K14").Value = "914,14 TL - 4.265,98 TL"

# Row 24
$ws.Range("D24").ClearContents()
$ws.Range("J24").ClearContents()

# Row 25
$ws.Range("D25").ClearContents()
$ws.Range("J25").ClearContents()
